$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the existing data so row 1 becomes a header
# row for the list (all existing rows shift down by one).
$ws.Rows("1:1").Insert() | Out-Null

# Give the new header row a title.
$ws.Range("A1").Value = "Product Description"

# Turn the whole list (header + all data rows) into a native Excel Table,
# using row 1 as the column header.
$rng = $ws.Range("A1:A42")
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $rng,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# Select the full table range, matching the resulting selection state.
$rng.Select() | Out-Null
